# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the affected rows, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.900.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.354.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.82%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.675'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.91%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.57'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.65%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +19.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.04'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.90%  '

$ws.Range("E12").Value = '  +2.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.703.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.65'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.901'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.358.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.889.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.65%  '

$ws.Range("E19").Value = '  +4.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '254.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.20%  '

$ws.Range("E25").Value = '  +3.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.16%  '

$ws.Range("E30").Value = '  +6.75%  '

$ws.Range("E31").Value = '  +1.86%  '

$ws.Range("E32").Value = '  +5.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0718'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.10%  '

$ws.Range("E37").Value = '  -1.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0267'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.72%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.81%  '

$ws.Range("E44").Value = '  +4.16%  '

$ws.Range("E45").Value = '  +1.06%  '

$ws.Range("E46").Value = '  +1.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.181'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.42%  '

$ws.Range("E49").Value = '  +4.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.435.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("E51").Value = '  +1.44%  '

